$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values could be misread as numbers
$textRefs = @("D5", "D6", "D8", "D9", "D11", "D13", "D14", "D19", "D20", "D21", "D22", "D23", "D24", "D28", "D29", "D32", "D35", "D36", "D38", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($ref in $textRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply cell value updates per row
# Row 2
$ws.Range("D2").Value = "68.317.35"
$ws.Range("E2").Value = "  +1.30%  "

# Row 3
$ws.Range("D3").Value = "3.562.70"
$ws.Range("E3").Value = "  +1.79%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "618.71"
$ws.Range("E5").Value = "  +2.07%  "

# Row 6
$ws.Range("D6").Value = "155.15"
$ws.Range("E6").Value = "  +4.01%  "

# Row 7
$ws.Range("D7").Value = "3.565.04"
$ws.Range("E7").Value = "  +1.90%  "

# Row 8
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.05%  "

# Row 9
$ws.Range("D9").Value = "0.491"
$ws.Range("E9").Value = "  +2.19%  "

# Row 10
$ws.Range("E10").Value = "  +5.23%  "

# Row 11
$ws.Range("D11").Value = "7.40"
$ws.Range("E11").Value = "  +6.64%  "

# Row 12
$ws.Range("E12").Value = "  +3.89%  "

# Row 13
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").Value = "33.24"
$ws.Range("E13").Value = "  +5.51%  "

# Row 14
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "0.0000221"
$ws.Range("E14").Value = "  +1.09%  "

# Row 15
$ws.Range("D15").Value = "4.166.03"
$ws.Range("E15").Value = "  +1.80%  "

# Row 16
$ws.Range("D16").Value = "3.574.23"
$ws.Range("E16").Value = "  +2.03%  "

# Row 17
$ws.Range("D17").Value = "68.354.23"
$ws.Range("E17").Value = "  +1.47%  "

# Row 18
$ws.Range("E18").Value = "  -0.13%  "

# Row 19
$ws.Range("D19").Value = "6.74"
$ws.Range("E19").Value = "  +5.72%  "

# Row 20
$ws.Range("D20").Value = "16.01"
$ws.Range("E20").Value = "  +6.76%  "

# Row 21
$ws.Range("D21").Value = "9.98"
$ws.Range("E21").Value = "  +11.43%  "

# Row 22
$ws.Range("D22").Value = "454.17"
$ws.Range("E22").Value = "  +1.83%  "

# Row 23
$ws.Range("D23").Value = "0.642"
$ws.Range("E23").Value = "  +3.84%  "

# Row 24
$ws.Range("D24").Value = "78.44"
$ws.Range("E24").Value = "  +1.48%  "

# Row 25
$ws.Range("E25").Value = "  +1.99%  "

# Row 26
$ws.Range("D26").Value = "3.707.20"
$ws.Range("E26").Value = "  +1.83%  "

# Row 27
$ws.Range("E27").Value = "  -0.24%  "

# Row 28
$ws.Range("D28").Value = "9.18"
$ws.Range("E28").Value = "  +11.45%  "

# Row 29
$ws.Range("D29").Value = "10.54"
$ws.Range("E29").Value = "  +3.57%  "

# Row 30
$ws.Range("E30").Value = "  +10.38%  "

# Row 31
$ws.Range("E31").Value = "  +3.20%  "

# Row 32
$ws.Range("D32").Value = "0.172"
$ws.Range("E32").Value = "  +4.22%  "

# Row 33
$ws.Range("E33").Value = "  +0.17%  "

# Row 34
$ws.Range("E34").Value = "  +3.53%  "

# Row 35
$ws.Range("D35").Value = "26.13"
$ws.Range("E35").Value = "  +1.91%  "

# Row 36
$ws.Range("D36").Value = "1.91"
$ws.Range("E36").Value = "  +3.67%  "

# Row 37
$ws.Range("D37").Value = "3.556.41"
$ws.Range("E37").Value = "  +1.93%  "

# Row 38
$ws.Range("D38").Value = "8.23"
$ws.Range("E38").Value = "  +3.27%  "

# Row 39
$ws.Range("E39").Value = "  +8.93%  "

# Row 40
$ws.Range("E40").Value = "  +0.04%  "

# Row 41
$ws.Range("D41").Value = "181.46"
$ws.Range("E41").Value = "  +4.64%  "

# Row 42
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.05%  "

# Row 43
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").Value = "0.0915"
$ws.Range("E43").Value = "  +4.66%  "

# Row 44
$ws.Range("E44").Value = "  +3.25%  "

# Row 45
$ws.Range("D45").Value = "31.11"
$ws.Range("E45").Value = "  +11.61%  "

# Row 46
$ws.Range("D46").Value = "0.898"
$ws.Range("E46").Value = "  +2.05%  "

# Row 47
$ws.Range("D47").Value = "46.18"
$ws.Range("E47").Value = "  +1.63%  "

# Row 48
$ws.Range("D48").Value = "1.33"
$ws.Range("E48").Value = "  +4.26%  "

# Row 49
$ws.Range("D49").Value = "2.66"
$ws.Range("E49").Value = "  +4.31%  "

# Row 50
$ws.Range("E50").Value = "  +3.50%  "

# Row 51
$ws.Range("D51").Value = "0.262"
$ws.Range("E51").Value = "  +7.73%  "

# Restore default style on cells where we forced text formatting
foreach ($ref in $textRefs) {
    $ws.Range($ref).Style = "Normal"
}
